$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C holds the Excel serial date 45188
# (2023-09-19) for every data row (rows 2-498). This update bumps that date
# by one day to 45189 (2023-09-20) across the whole column.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 498) { $lastRow = 498 }

$ws.Range("C2:C$lastRow").Value = 45189
